$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 10; $row++) {
    $ws.Range("C$row").Value = 45207
}
